$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new wishers to the existing "Wishes" lists so that each mutual
# group also references someone outside of its own mutual pool.
$ws.Range("C44").Value = "Pertti Koivisto, Laura Koivisto, Elisabet Koivisto, Reijo Koivisto, Satu Räsänen"
$ws.Range("C47").Value = "Tuomas Räsänen, Mari Räsänen, Olli Räsänen, Heidi Räsänen, Martti Koivisto"
$ws.Range("C132").Value = "Kati Kuusisto, Julius Kuusisto, Margareta Kuusisto, Topi Kuusisto, Elli Tamminen"
$ws.Range("C136").Value = "Kim Kuusisto, Kati Kuusisto, Julius Kuusisto, Margareta Kuusisto, Alina Tamminen"
$ws.Range("C139").Value = "Alina Tamminen, Johan Tamminen, Konsta Tamminen, Marja-Liisa Tamminen, Kim Kuusisto, toimisko tämä?"

# Leave the view where the user ended up after the edit (last cell touched).
$ws.Range("C139").Select()
